$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new activity-log entry in row 10:
#   last 4 digits | date | starttime | endtime | description
$ws.Range("B10").Value = 6977
$ws.Range("C10").Value = 43923
$ws.Range("D10").Value = 0.3430555555555555
$ws.Range("E10").Value = 0.3444444444444445
$ws.Range("G10").Value = "Updated Final Project PDF document with own name and student number."

# Move the active selection to A10, matching the author's cursor position
# when they saved the workbook.
$ws.Range("A10").Select()
